$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "polite_expressions" value in row 10 (C10), leaving it blank.
$ws.Range("C10").Value = ""

# Add a new row 11 duplicating the "id"/"source_file"/"text" info that
# originally lived in row 8, paired with the annotator/score/issue pattern
# used by row 10.
$ws.Range("A11").Value = "parisk"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "nan"
$ws.Range("D11").Value = "APC"
$ws.Range("E11").Value = "RES"
$ws.Range("F11").Value = "d3fb2dcb-ee08-4432-9f4b-c252dbb3433f"
$ws.Range("G11").Value = "SJ3dBGZ0Z_annotated.xlsx"
$ws.Range("H11").Value = "We evaluate our method on NLP task for two reasons: 1) they are particularly well-suited for evaluating our method (naturally large output spaces) 2) we did not dispose of the computational resources to tackle tasks from other domains such as vision (e.g. Flickr100M) which requires hundreds of GPUs for weeks."
$ws.Range("I11").Value = "Correct"

Write-Host "Edit applied"
